$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = 1.012508749961853
$ws.Range("B1").Value = 1.869760036468506
$ws.Range("C1").Value = 7.718123435974121
$ws.Range("D1").Value = 2.520470142364502
$ws.Range("E1").Value = 0.4364215135574341
